$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-03-30 Sunday"

# Update the multiplication problems in the table, addressed by (row, col) and
# assigning Range.Text directly (rather than Find/Replace) so that cells whose old
# or new values collide with other cells in the table (e.g. 684x7= -> 733x9= while
# another cell already has 733x9= -> 561x8=) cannot cross-contaminate each other.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "616×7="
$t.Cell(1,2).Range.Text = "909×8="
$t.Cell(1,3).Range.Text = "215×8="
$t.Cell(1,4).Range.Text = "733×9="
$t.Cell(1,5).Range.Text = "875×7="
$t.Cell(5,1).Range.Text = "747×8="
$t.Cell(5,2).Range.Text = "334×9="
$t.Cell(5,3).Range.Text = "695×9="
$t.Cell(5,4).Range.Text = "863×2="
$t.Cell(5,5).Range.Text = "195×8="
$t.Cell(10,1).Range.Text = "643×4="
$t.Cell(10,2).Range.Text = "129×7="
$t.Cell(10,3).Range.Text = "387×8="
$t.Cell(10,4).Range.Text = "449×3="
$t.Cell(10,5).Range.Text = "146×3="
$t.Cell(15,1).Range.Text = "319×6="
$t.Cell(15,2).Range.Text = "348×6="
$t.Cell(15,3).Range.Text = "584×8="
$t.Cell(15,4).Range.Text = "122×4="
$t.Cell(15,5).Range.Text = "400×2="
$t.Cell(20,1).Range.Text = "265×8="
$t.Cell(20,2).Range.Text = "959×5="
$t.Cell(20,3).Range.Text = "339×3="
$t.Cell(20,4).Range.Text = "561×8="
$t.Cell(20,5).Range.Text = "410×5="
